$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("summary")
$ws.Range("A2").Value = 2904.199958328551
$ws.Range("D2").Value = 51.70103124997319

$ws = $wb.Worksheets.Item("bus")
$ws.Range("B2").Value = -6.954336164209468
$ws.Range("C2").Value = 1.041171740726286
$ws.Range("B3").Value = -6.987038895179992
$ws.Range("C3").Value = 1.035000000900278
$ws.Range("B4").Value = -5.791314030672015
$ws.Range("C4").Value = 1.049999989553545
$ws.Range("B5").Value = -9.380426072626054
$ws.Range("C5").Value = 1.006629803587966
$ws.Range("B6").Value = -9.703135100475894
$ws.Range("C6").Value = 1.023816162635372
$ws.Range("B7").Value = -12.18984785700471
$ws.Range("C7").Value = 1.015735577424695
$ws.Range("B8").Value = -7.0138043603139
$ws.Range("C8").Value = 1.025000000861727
$ws.Range("B9").Value = -10.7918773526312
$ws.Range("C9").Value = 0.996994344779235
$ws.Range("B10").Value = -7.255529263787974
$ws.Range("C10").Value = 1.016761895503189
$ws.Range("B11").Value = -9.343244474360933
$ws.Range("C11").Value = 1.032710967342917
$ws.Range("B12").Value = -2.175645845717776
$ws.Range("C12").Value = 0.993463264623684
$ws.Range("B13").Value = -1.519622106706469
$ws.Range("C13").Value = 1.007024443039843
$ws.Range("B14").Value = 0.0000000000000000000001922255097482938
$ws.Range("C14").Value = 1.020000000687559
$ws.Range("B15").Value = 2.116497042707441
$ws.Range("C15").Value = 0.9800000029333509
$ws.Range("B16").Value = 10.47761840216531
$ws.Range("C16").Value = 1.09578684298166
$ws.Range("B17").Value = 9.901490906490791
$ws.Range("C17").Value = 1.03552312332949
$ws.Range("B18").Value = 14.21567334744549
$ws.Range("C18").Value = 1.044723716059756
$ws.Range("B19").Value = 15.53637998996704
$ws.Range("C19").Value = 1.049999999808577
$ws.Range("B20").Value = 8.557752832219609
$ws.Range("C20").Value = 1.033921973285685
$ws.Range("B21").Value = 9.284541911580224
$ws.Range("C21").Value = 1.042247791338515
$ws.Range("B22").Value = 16.29525361727388
$ws.Range("C22").Value = 1.050000356255982
$ws.Range("B23").Value = 21.98422090266471
$ws.Range("C23").Value = 1.050000000112305
$ws.Range("B24").Value = 10.38975485847643
$ws.Range("C24").Value = 1.050000000393651
$ws.Range("B25").Value = 4.964784838117613
$ws.Range("C25").Value = 1.027050270768929

$ws = $wb.Worksheets.Item("generator")
$ws.Range("E2").Value = 10.469
$ws.Range("E3").Value = 10.469
$ws.Range("E4").Value = 10.469
$ws.Range("E5").Value = 10.469
$ws.Range("E6").Value = -9.765000000000001
$ws.Range("E7").Value = -9.765000000000001
$ws.Range("E8").Value = -9.765000000000001
$ws.Range("E9").Value = -9.765000000000001
$ws.Range("E17").Value = 137.037
$ws.Range("E18").Value = 137.037
$ws.Range("E19").Value = 137.037
$ws.Range("E20").Value = 137.037
$ws.Range("E21").Value = 137.037
$ws.Range("E22").Value = 137.037
$ws.Range("E23").Value = -204.612
$ws.Range("E24").Value = 93.45399999999999
$ws.Range("E25").Value = -239.835

$ws = $wb.Worksheets.Item("branch")
$ws.Range("D2").Value = 231.0974088300875
$ws.Range("E2").Value = -224.9052727344524
$ws.Range("F2").Value = 6.192136095635004
$ws.Range("D3").Value = -48.45925443476258
$ws.Range("E3").Value = 49.60070844371397
$ws.Range("F3").Value = 1.141454008951392
$ws.Range("D4").Value = 177.8623638468168
$ws.Range("E4").Value = -177.3404277294767
$ws.Range("F4").Value = 0.5219361173400383
$ws.Range("D5").Value = 159.5927872348254
$ws.Range("E5").Value = -157.5706497575499
$ws.Range("F5").Value = 2.022137477275576
$ws.Range("D6").Value = 369.3602581928429
$ws.Range("E6").Value = -362.277622993714
$ws.Range("F6").Value = 7.082635199128884
$ws.Range("D7").Value = -112.9343918605495
$ws.Range("E7").Value = 114.9999970027528
$ws.Range("F7").Value = 2.065605142203308
$ws.Range("D8").Value = -28.30292709544883
$ws.Range("E8").Value = 28.67024434418971
$ws.Range("F8").Value = 0.3673172487408782
$ws.Range("D9").Value = 315.3324155686497
$ws.Range("E9").Value = -312.323679657658
$ws.Range("F9").Value = 3.008735910991645
$ws.Range("D10").Value = 36.87525852921386
$ws.Range("E10").Value = -36.52244851399961
$ws.Range("F10").Value = 0.3528100152142499
$ws.Range("D11").Value = -220.9042426163613
$ws.Range("E11").Value = 224.6073084525962
$ws.Range("F11").Value = 3.703065836234876
$ws.Range("D12").Value = 168.2776219955612
$ws.Range("E12").Value = -166.5624347190722
$ws.Range("F12").Value = 1.715187276489027
$ws.Range("D13").Value = 86.87483807478857
$ws.Range("E13").Value = -86.30134756302856
$ws.Range("F13").Value = 0.5734905117600109
$ws.Range("D14").Value = 226.8349498303853
$ws.Range("E14").Value = -221.5777937823791
$ws.Range("F14").Value = 5.257156048006184
$ws.Range("D15").Value = 101.0338173711477
$ws.Range("E15").Value = -100.7580177567392
$ws.Range("F15").Value = 0.2757996144084673
$ws.Range("D16").Value = -110.0028256249187
$ws.Range("E16").Value = 112.7129099135157
$ws.Range("F16").Value = 2.710084288597092
$ws.Range("D17").Value = 140.4072067707107
$ws.Range("E17").Value = -137.991987839173
$ws.Range("F17").Value = 2.41521893153771
$ws.Range("D18").Value = 10.93761504779571
$ws.Range("E18").Value = -10.87699104974649
$ws.Range("F18").Value = 0.06062399804922375
$ws.Range("D19").Value = 55.52480781459366
$ws.Range("E19").Value = -55.43118392294721
$ws.Range("F19").Value = 0.0936238916464438
$ws.Range("D20").Value = 38.732361113112
$ws.Range("E20").Value = -38.10815641637365
$ws.Range("F20").Value = 0.6242046967383508
$ws.Range("D21").Value = 223.2605150646424
$ws.Range("E21").Value = -219.6601124552858
$ws.Range("F21").Value = 3.600402609356612
$ws.Range("D22").Value = 88.57846656106835
$ws.Range("E22").Value = -87.54074556523743
$ws.Range("F22").Value = 1.037720995830915
$ws.Range("D23").Value = -37.47755148600039
$ws.Range("E23").Value = 37.97975438722105
$ws.Range("F23").Value = 0.5022029012206608
$ws.Range("D24").Value = -12.58046682726585
$ws.Range("E24").Value = 12.63044818795772
$ws.Range("F24").Value = 0.04998136069186976
$ws.Range("D25").Value = 36.75801775673923
$ws.Range("E25").Value = -36.68037120899265
$ws.Range("F25").Value = 0.07764654774657376
$ws.Range("D26").Value = 101.0338173711477
$ws.Range("E26").Value = -100.7580177567392
$ws.Range("F26").Value = 0.2757996144084673
$ws.Range("D27").Value = -61.45811387440724
$ws.Range("E27").Value = 62.24653886545782
$ws.Range("F27").Value = 0.7884249910505847
$ws.Range("D28").Value = -107.6392575820147
$ws.Range("E28").Value = 107.9662455406568
$ws.Range("F28").Value = 0.326987958642122
$ws.Range("D29").Value = 20.2413939355046
$ws.Range("E29").Value = -19.95745172307686
$ws.Range("F29").Value = 0.283942212427743
$ws.Range("D30").Value = 59.90295271036039
$ws.Range("E30").Value = -59.66698325979712
$ws.Range("F30").Value = 0.235969450563267
$ws.Range("D31").Value = 9.575011987546883
$ws.Range("E31").Value = -9.541886125592757
$ws.Range("F31").Value = 0.03312586195412526
$ws.Range("D32").Value = 36.75801775673923
$ws.Range("E32").Value = -36.68037120899265
$ws.Range("F32").Value = 0.07764654774657376
$ws.Range("D33").Value = 223.2605150646424
$ws.Range("E33").Value = -219.6601124552858
$ws.Range("F33").Value = 3.600402609356612
$ws.Range("D34").Value = 55.52480781459366
$ws.Range("E34").Value = -55.43118392294721
$ws.Range("F34").Value = 0.0936238916464438

$ws = $wb.Worksheets.Item("transformer")
$ws.Range("D2").Value = 103.751360831875
$ws.Range("E2").Value = -103.498528788644
$ws.Range("F2").Value = 0.2528320432310105
$ws.Range("D3").Value = 119.1440728202233
$ws.Range("E3").Value = -118.8061637582331
$ws.Range("F3").Value = 0.3379090619902048
$ws.Range("D4").Value = 149.1124214502258
$ws.Range("E4").Value = -148.5873693891832
$ws.Range("F4").Value = 0.5250520610425191
$ws.Range("D5").Value = 220.9042426163613
$ws.Range("E5").Value = -219.6078593919854
$ws.Range("F5").Value = 1.296383224375885
$ws.Range("G5").Value = 1.100000010873496
$ws.Range("D6").Value = 165.4281831740263
$ws.Range("E6").Value = -164.8075030949366
$ws.Range("F6").Value = 0.6206800790897082

